$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (pushes existing row 14 and below down by one),
# inheriting formatting from the row above.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row with the new docente responsável entry.
$ws.Range("B14").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("C14").Value = "1341653 - Maria José Ramos Sandim"
